$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "Conditional offer 이신 분들한테 연락드립니다"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/notice-to-conditional-offers/#utm_source=rss&utm_medium=rss&utm_campaign=notice-to-conditional-offers"

$ws.Range("D28").Value = "[임피던스 제어] Direct Impedance Modulation"
$ws.Range("E28").Value = "https://ropiens.tistory.com/148"

$ws.Range("D29").Value = "[논문리뷰] Adversarial Latent Autoencoders"
$ws.Range("E29").Value = "https://blog.promedius.ai/alae_1/"
